# Settings.xlsx edit: update simulation "Values" row parameters and
# refresh the current cell selection on the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 ("Values") parameter updates
$ws.Range("C4").Value = 10          # Spring constant (k)
$ws.Range("D4").Value = 0.01        # Surface tension coefficient (sigma)
$ws.Range("E4").Value = 0.05        # Damping coefficient (c)
$ws.Range("I4").Value = "150"       # Repeats (stored as text)

# Move the active selection to F12
$ws.Range("F12").Select()
